# "Generate Report for Handoff"
#
# The localization-status report is regenerated: every locale that was
# previously "Handed back: in sync with en-US" is now "Ready for handoff"
# (a fresh handoff package was produced), and the associated timestamps
# on the Overview sheet and the per-locale (zh-cn / de-de) sheets move
# forward a few seconds to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

# ----- Overview sheet --------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "Ready for handoff"          # zh-cn status
$ovw.Range("F2").Value = "Ready for handoff"          # de-de status
$ovw.Range("G2").Value = "2016-09-05 17:08:58"        # Latest HO Xliff Generate Date

# ----- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"         # Status
$zhcn.Range("H2").Value = "2016-09-05 17:08:54"       # Latest Handoff Datetime

# ----- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"         # Status
$dede.Range("H2").Value = "2016-09-05 17:08:58"       # Latest Handback DateTime (shared w/ Overview G2)

# ----- Column widths -----------------------------------------------------
# The Status text shrank ("Handed back: in sync with en-US" -> "Ready for
# handoff"), so the column that displays it re-autofits narrower on both
# the Overview sheet (zh-cn/de-de status columns E & F) and on each
# per-locale sheet (Status column C).
$ovw.Columns.Item(5).ColumnWidth = 16.35
$ovw.Columns.Item(6).ColumnWidth = 16.35
$zhcn.Columns.Item(3).ColumnWidth = 16.35
$dede.Columns.Item(3).ColumnWidth = 16.35
